# Generate Report for Handback
# Updates the localization-status report after a successful handback:
#  - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#    on the Overview sheet and on each per-language sheet.
#  - The zh-cn / de-de "Latest Handback DateTime" gets refreshed to the new
#    handback timestamp.
#  - The (now resolved) "Error Detail" stale-handback warning is cleared.

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

# ---- Overview sheet -------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusNew
$wsOverview.Range("F2").Value = $statusNew

# Widen the zh-cn / de-de status columns to fit the longer text.
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---- zh-cn sheet ------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $statusNew
$wsZhCn.Range("K2").Value = "2016-09-02 20:57:29"
$wsZhCn.Range("P2").ClearContents() | Out-Null

$wsZhCn.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZhCn.Columns.Item(16).ColumnWidth = 12.833333333333334

# ---- de-de sheet --------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $statusNew
$wsDeDe.Range("K2").Value = "2016-09-02 20:57:37"
$wsDeDe.Range("P2").ClearContents() | Out-Null

$wsDeDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDeDe.Columns.Item(16).ColumnWidth = 12.833333333333334
